# Minor correction to Powerpoint
#
# 1) The "Updated automatically" date field cached on the Slide Master and
#    on every Slide Layout shows a stale date (3/20/17). Refresh the
#    cached text to 4/10/17 on each one.
# 2) Slide 19's "TextBox 24" heading was mis-numbered ("7: ADVERSE
#    REACTION...") and should read "8: ADVERSE REACTION...".

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "3/20/17") {
                $tr.Text = "4/10/17"
            }
        }
    }
}

# Slide Master footer date field.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every Slide Layout's footer date field.
for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# Slide 19 heading textbox: "7" -> "8", keeping the rest of the
# run-formatting (bold/underline/italic) untouched.
$slide19 = $p.Slides.Item(19)
for ($i = 1; $i -le $slide19.Shapes.Count; $i++) {
    $shp = $slide19.Shapes.Item($i)
    if ($shp.Name -eq "TextBox 24") {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "7: ADVERSE REACTION: Check reason(s) for regimen change") {
            $firstChar = $tr.Characters(1, 1)
            $firstChar.Text = "8"

            # Re-apply the underline on the colon+space that follows the
            # number so it becomes its own run, matching how the author's
            # edit split "7: ADVERSE REACTION" into "8" + ": " + "ADVERSE REACTION".
            $colonSpace = $tr.Characters(2, 2)
            $colonSpace.Font.Underline = $true
        }
    }
}
